$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete original row 4 (even_MAG-GUT27169.fa)
$ws.Rows.Item(4).Delete()

# After the above delete, original row 6 (even_MAG-GUT38735.fa) is now row 5
$ws.Rows.Item(5).Delete()
